# Fill in the four empty "date" cells (first column of the date pair) for
# the ПР09 / Л07 / ПР10 / ЛР07 rows of the schedule table with their
# corresponding dates: 09.05, 10.05, 12.05, 13.05.
#
# Each inserted run must carry:
#   rFonts ascii/hAnsi/cs = "Times New Roman", sz = 28, szCs = 28
#
# Because this runtime treats a Range handle as potentially "stale" after
# a structural edit (Text assignment / Find.Execute both count), every
# object (table/row/cell/range) is re-fetched fresh from
# $word.ActiveDocument right before each step that mutates formatting.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Row index (in the single schedule table) -> date text to insert into
# that row's first ("day.month", 1133-twip-wide) cell, which is empty.
$targets = @(
    @{ Row = 24; Date = "09.05" },  # ПР09 row
    @{ Row = 25; Date = "10.05" },  # Л07 row
    @{ Row = 26; Date = "12.05" },  # ПР10 row
    @{ Row = 27; Date = "13.05" }   # ЛР07 row
)

foreach ($item in $targets) {
    $rowIdx = $item.Row
    $dateText = $item.Date

    # 1) Put the plain text into the (currently empty) cell.
    $d1 = $word.ActiveDocument
    $t1 = $d1.Tables.Item(1)
    $cell1 = $t1.Rows.Item($rowIdx).Cells.Item(1)
    $cell1.Range.Text = $dateText

    # 2) Re-fetch and apply the ascii/hAnsi + complex-script font name via
    #    a range-scoped Find/Replace (only way this runtime will stamp the
    #    w:cs="Times New Roman" attribute onto the run's rFonts).
    $d2 = $word.ActiveDocument
    $t2 = $d2.Tables.Item(1)
    $r2 = $t2.Rows.Item($rowIdx).Cells.Item(1).Range
    $find = $r2.Find
    $find.ClearFormatting()
    $find.Replacement.ClearFormatting()
    $find.Replacement.Font.Name = "Times New Roman"
    $find.Replacement.Font.NameBi = "Times New Roman"
    $find.Execute($dateText, $false, $false, $false, $false, $false, $true, 0, $false, $dateText, 2)

    # 3) Re-fetch again and set the sizes (28 half-points => 14pt) for both
    #    the Western and complex-script runs (sz / szCs).
    $d3 = $word.ActiveDocument
    $t3 = $d3.Tables.Item(1)
    $r3 = $t3.Rows.Item($rowIdx).Cells.Item(1).Range
    $r3.Font.Size = 14
    $r3.Font.SizeBi = 14
}
